$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: change E3 formula and add F3/G3 values ---
$ws.Range("E3").Formula = "=SUM(F3:G3)"
$ws.Range("E3").Font.Bold = $true
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3

# --- Row 4: add new "Streda 24.4" day column with its total and daily hours ---
$ws.Range("D4").Value = "Streda 24.4"
$ws.Range("E4").Formula = "=SUM(F4:J4)"
$ws.Range("E4").Font.Bold = $true
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 1

# --- Row 13: add weekly total for column E ---
$ws.Range("E13").Formula = "=SUM(E2:E12)"

# --- Update selection to match the final state ---
$ws.Range("H12").Select()

$wb.Save()
